# Update BunkerPrices at 2025-04-14 03:01
# Adds a new data row (row 30) to Sheet1 and updates the number format
# of the "Date" cell in the previous last row (Y29) from date-only to
# date+time, matching the style used for the newly appended row's date
# cell moving forward.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29's Date cell (Y29) switches from the date-only style to the
# date+time style.
$ws.Range("Y29").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row of bunker price data, in column order A..AV.
$newRowValues = @(
    567, 478, 442, 549, 494, 532, 475, 569, 512, 442, 571, 483, 447, 505,
    555, 483, 619, 495, 475, 480, 620, 530, 589, 475, 45757, 846, 555,
    519.5, 512, 538, 500, 502, 750, 459, 740, 475, 486, 550, 535, 487,
    536, 523, 563, 545, 620, 632, 489, 470
)

$rowIndex = 30
for ($col = 1; $col -le $newRowValues.Length; $col++) {
    $ws.Cells.Item($rowIndex, $col).Value = $newRowValues[$col - 1]
}

# Column Y (25) holds the date; it keeps the date-only style that Y29
# previously had.
$ws.Cells.Item($rowIndex, 25).NumberFormat = "YYYY-MM-DD"
